$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("POWER_CONTROL")

$ports = @("port 1","port 2","port 3","port 4","port 5","port 6","port 7","port 8")
$remarks = @("1 - OFF, 2 - ON","2 - OFF, 2 - ON","3 - OFF, 2 - ON","4 - OFF, 2 - ON","5 - OFF, 2 - ON","6 - OFF, 2 - ON","7 - OFF, 2 - ON","8 - OFF, 2 - ON")

# 1) Populate column C values first (ports) so the shared strings land 260..267
for ($i = 0; $i -lt 8; $i++) {
  $row = $i + 2
  $ws.Range("C$row").Value = $ports[$i]
}

# 2) Populate column D values + number format + centering (remarks land 268..275,
#    and this also creates the new numFmt164+center style)
for ($i = 0; $i -lt 8; $i++) {
  $row = $i + 2
  $dcell = $ws.Range("D$row")
  $dcell.Value = $remarks[$i]
  $dcell.NumberFormat = "0.0"
  $dcell.HorizontalAlignment = -4108
}

# 3) Center column C (creates the plain center-only style)
for ($i = 0; $i -lt 8; $i++) {
  $row = $i + 2
  $ws.Range("C$row").HorizontalAlignment = -4108
}

# Header D1 alignment (style 3 bold -> style 1 bold+center)
$ws.Range("D1").HorizontalAlignment = -4108
